# Append the 2025-06-10 resazurin mortality readings (samples A1-E5) as a
# new 25-row block directly below the existing 2025-06-03 block, mirroring
# the same structure/formatting used for every prior date block, then leave
# the sheet scrolled/selected the way the author left it (B146 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 2025-06-03 block (rows 102-126) is the most recent complete 25-row
# A1..E5 block. Copy just columns A:C (skip D, which only holds a couple of
# one-off "remove sample" notes) into the 25 new rows starting at A127 so the
# new block inherits the same number formatting / styles (s="1" on B & C).
$ws.Range("A102:C126").Copy($ws.Range("A127"))

# Stamp the new block with the new collection date.
$ws.Range("A127:A151").Value2 = 20250610

# The source block had two non-zero / "NA" mortality notes (from samples B1
# and C3) that don't apply to this new reading - reset every mortality cell
# in the new block back to a plain 0, matching the rest of the data set.
$ws.Range("C127:C151").Value2 = 0

# Match the author's final scroll position / active cell in the sheet.
[void]$ws.Range("B146").Select()
